# Larson recommended changes to Model dev calendar
#
# The "Model Development" section used to contain a single blank spacer
# row (row 17) before the header. It's being expanded into five new
# sub-tasks (rows 17-21), which pushes every row below down by four.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert four fresh rows right after the existing blank row 17 - this
# shifts the old "Model Development" header (and everything after it)
# from row 18 down to row 22, while row 17 itself is reused for the
# first new sub-task.
$ws.Rows("18:21").Insert()

# Row 17: Binary Sidewalk detector
$ws.Range("B17").Value = "Binary Sidewalk detector, investigate convolutional architectures for automatically labeling sidewalks in images (10 days)"
$ws.Range("C17").Value = 42903
$ws.Range("D17").Value = 10
$ws.Range("E17").Formula = "=C17+D17"

# Row 18: Sidewalk quality detection
$ws.Range("B18").Value = "Sidewalk quality detection, adapt architectures for detection and quality classification"
$ws.Range("C18").Value = 42911
$ws.Range("D18").Value = 15
$ws.Range("E18").Formula = "=C18+D18"

# Row 19: Parameters tuning and further architecture exploration
$ws.Range("B19").Value = "Parameters tuning and further architecture exploration"
$ws.Range("C19").Value = 42921
$ws.Range("D19").Value = 15
$ws.Range("E19").Formula = "=C19+D19"

# Row 20: Evaluation of models to other neighborhoods
$ws.Range("B20").Value = "Evaluation of models to other neighborhoods"
$ws.Range("C20").Value = 42931
$ws.Range("D20").Value = 10
$ws.Range("E20").Formula = "=C20+D20"

# Row 21: Identification of individual labels that contribute to overall rating
$ws.Range("B21").Value = " Identification of individual labels that contribute to overall rating"
$ws.Range("C21").Value = 42941
$ws.Range("D21").Value = 20
$ws.Range("E21").Formula = "=C21+D21"

# New row 17 (the first of the group) gets a taller custom row height to
# fit its long wrapped label; row 5's wrap also re-flows slightly.
$ws.Rows(5).RowHeight = 25.5
$ws.Rows(17).RowHeight = 39

# Point the chart's two series at the now-larger data range
# (Sheet1!A2:B32 / C2:C32 / D2:D32 instead of ...:28).
$co = $ws.ChartObjects(1)
$chart = $co.Chart
$s1 = $chart.SeriesCollection(1)
$s1.Formula = "=SERIES(Sheet1!`$C`$1,Sheet1!`$A`$2:`$B`$32,Sheet1!`$C`$2:`$C`$32,1)"
$s2 = $chart.SeriesCollection(2)
$s2.Formula = "=SERIES(Sheet1!`$D`$1,Sheet1!`$A`$2:`$B`$32,Sheet1!`$D`$2:`$D`$32,2)"

# The chart grew taller/shifted up slightly to keep pace with the extra
# rows - reposition its anchor to match.
$co.Top = 48.75
$co.Height = 612

# Reset the view: scroll back to column A (drop the old topLeftCell="E1")
# and leave the selection on C18, matching where the editor was working.
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("C18").Select()
